$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (these strings are shared across rows 2:5) ---

# The date column (A2:A5) holds a plain text string "2025-12-03" that must
# become "2025-12-05". A bare string assignment gets auto-recognized by
# Excel as a date (and reformatted to a date serial), so force it in as
# literal text via a leading apostrophe, then clear the resulting
# quote-prefix/format so the cell ends up with no explicit style, matching
# how the original cell was stored.
$dateRange = $ws.Range("A2:A5")
$dateRange.Value = "'2025-12-05"
$dateRange.ClearFormats()

# The judgment column (O2:O5) just holds plain text, no special handling needed.
$ws.Range("O2:O5").Value = "⚪ 중립 구간"

# --- Row 2: HD HYUNDAI MIPO ---
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 55.7
$ws.Range("N2").Value = 52.43913937059539

# --- Row 3: HDKSOE ---
$ws.Range("D3").Value = 419000
$ws.Range("E3").Value = 45.3
$ws.Range("F3").Value = -1.87
$ws.Range("H3").Value = 70
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 48.1
$ws.Range("N3").Value = 52.43913937059539

# --- Row 4: SamsungHvyInd ---
$ws.Range("D4").Value = 24650
$ws.Range("E4").Value = 36.4
$ws.Range("F4").Value = -2.57
$ws.Range("H4").Value = 56
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 47.9
$ws.Range("N4").Value = 52.43913937059539

# --- Row 5: Hanwha Ocean ---
$ws.Range("D5").Value = 106500
$ws.Range("E5").Value = 20.6
$ws.Range("F5").Value = -3.62
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 46.7
$ws.Range("N5").Value = 52.43913937059539

$wb.Save()
